$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.516.16"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "1.602.21"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0597"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0913"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "1.831.36"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "1.603.13"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "29.511.32"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.535"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("E28").Value = "  +4.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("D34").Value = "1.428.54"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.95%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.534"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "54.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +24.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.794"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "1.741.77"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  -3.78%  "
